$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.410.90"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.849.51"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.77"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6301"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07689"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07747"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "1.846.75"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.025"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001084"
$ws.Range("E14").Value = "  +8.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6808"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.75"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "2.104.69"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.150"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "29.423.44"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.25"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.455"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.31"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.317"
$ws.Range("E29").Value = "  +4.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.469"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05724"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.054"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.849"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7089"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.779"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01797"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").Value = "1.224.56"
$ws.Range("E40").Value = "  -2.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.447"
$ws.Range("E41").Value = "  +4.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9109"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "2.013.47"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.87"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.18"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.139"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.042"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.689"
$ws.Range("E51").Value = "  +0.22%  "
